# Auto-generated edit script: update market price columns (H-N) per scheduled runner pull
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 523.1111
$ws.Range("I28").Value = 523.1111
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 523.1111
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -38.11109999999996
$ws.Range("N28").ClearContents()

$ws.Range("H62").Value = 7809.6
$ws.Range("I62").Value = 9002.5
$ws.Range("J62").Value = 7014.3335
$ws.Range("K62").Value = 9002.5
$ws.Range("L62").Value = 7014.3335
$ws.Range("M62").Value = -8378.5
$ws.Range("N62").Value = -8262.333500000001

$ws.Range("H65").Value = 7809.6
$ws.Range("I65").Value = 9002.5
$ws.Range("J65").Value = 7014.3335
$ws.Range("K65").Value = 45012.5
$ws.Range("L65").Value = 35071.6675
$ws.Range("M65").Value = -41892.5
$ws.Range("N65").Value = -41311.6675

$ws.Range("H112").Value = 3070.2727
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 3327.3
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 9981.900000000001
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -12197.9

$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -10008

$ws.Range("H129").Value = 4556.6
$ws.Range("I129").Value = 4556.6
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 13669.8
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -8669.800000000001

$ws.Range("H132").Value = 4776.5415
$ws.Range("I132").Value = 1379.0834
$ws.Range("J132").Value = 8174
$ws.Range("K132").Value = 4137.2502
$ws.Range("L132").Value = 24522
$ws.Range("M132").Value = -1607.2502
$ws.Range("N132").Value = -29582

$ws.Range("H137").Value = 2560.577
$ws.Range("I137").Value = 1824.909
$ws.Range("J137").Value = 3100.0667
$ws.Range("K137").Value = 5474.727000000001
$ws.Range("L137").Value = 9300.2001
$ws.Range("M137").Value = -2924.727000000001
$ws.Range("N137").Value = -14400.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 424.25
$ws.Range("I5").Value = 424.16666
$ws.Range("J5").Value = 424.5
$ws.Range("K5").Value = 424.16666
$ws.Range("L5").Value = 424.5
$ws.Range("M5").Value = -312.16666
$ws.Range("N5").Value = -648.5

$ws.Range("H109").Value = 99000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 99000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 99000
$ws.Range("N109").Value = -101774

$ws.Range("H122").Value = 2213.2856
$ws.Range("I122").Value = 2188.8
$ws.Range("J122").Value = 2274.5
$ws.Range("K122").Value = 6566.400000000001
$ws.Range("L122").Value = 6823.5
$ws.Range("M122").Value = -4116.400000000001
$ws.Range("N122").Value = -11723.5

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 424.25
$ws.Range("I4").Value = 424.16666
$ws.Range("J4").Value = 424.5
$ws.Range("K4").Value = 424.16666
$ws.Range("L4").Value = 424.5
$ws.Range("M4").Value = -309.16666
$ws.Range("N4").Value = -654.5

$ws.Range("H105").Value = 5296.1
$ws.Range("I105").Value = 4678
$ws.Range("J105").Value = 5708.1665
$ws.Range("K105").Value = 4678
$ws.Range("L105").Value = 5708.1665
$ws.Range("M105").Value = -2931
$ws.Range("N105").Value = -9202.166499999999

$ws.Range("H111").Value = 40000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 40000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5623.636
$ws.Range("I31").Value = 3692.5715
$ws.Range("J31").Value = 9003
$ws.Range("K31").Value = 3692.5715
$ws.Range("L31").Value = 9003
$ws.Range("M31").Value = -3397.5715
$ws.Range("N31").Value = -9593

$ws.Range("H34").Value = 5623.636
$ws.Range("I34").Value = 3692.5715
$ws.Range("J34").Value = 9003
$ws.Range("K34").Value = 3692.5715
$ws.Range("L34").Value = 9003
$ws.Range("M34").Value = -3490.5715
$ws.Range("N34").Value = -9407

$ws.Range("H99").Value = 5381.636
$ws.Range("I99").Value = 4465.722
$ws.Range("J99").Value = 9503.25
$ws.Range("K99").Value = 4465.722
$ws.Range("L99").Value = 9503.25
$ws.Range("M99").Value = -2967.722
$ws.Range("N99").Value = -12499.25

$ws.Range("H107").Value = 910.875
$ws.Range("I107").Value = 898.2857
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 898.2857
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1021.7143
$ws.Range("N107").Value = -4839

$ws.Range("H126").Value = 5381.636
$ws.Range("I126").Value = 4465.722
$ws.Range("J126").Value = 9503.25
$ws.Range("K126").Value = 13397.166
$ws.Range("L126").Value = 28509.75
$ws.Range("M126").Value = -10927.166
$ws.Range("N126").Value = -33449.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H131").Value = 1439.4286
$ws.Range("I131").Value = 1294.3334
$ws.Range("J131").Value = 1548.25
$ws.Range("K131").Value = 3883.0002
$ws.Range("L131").Value = 4644.75
$ws.Range("M131").Value = 1156.9998
$ws.Range("N131").Value = -14724.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3504.75
$ws.Range("I102").Value = 2009.5
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2009.5
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -387.5
$ws.Range("N102").Value = -8244

$ws.Range("H122").Value = 2470.5293
$ws.Range("I122").Value = 2635.6428
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 7906.928400000001
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -5456.928400000001
$ws.Range("N122").Value = -10000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5373.25
$ws.Range("I7").Value = 4996.5
$ws.Range("J7").Value = 5750
$ws.Range("K7").Value = 4996.5
$ws.Range("L7").Value = 5750
$ws.Range("M7").Value = -4884.5
$ws.Range("N7").Value = -5974

$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 900
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -730

$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1250
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -955
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 750
$ws.Range("I27").Value = 1250
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 1250
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -1143
$ws.Range("N27").Value = -714

$ws.Range("H36").Value = 39998
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 39998
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 39998
$ws.Range("N36").Value = -41122

$ws.Range("H40").Value = 3294.4
$ws.Range("I40").Value = 2618
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 2618
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -2482
$ws.Range("N40").Value = -6272

$ws.Range("H55").Value = 866.1905
$ws.Range("I55").Value = 265.66666
$ws.Range("J55").Value = 1106.4
$ws.Range("K55").Value = 265.66666
$ws.Range("L55").Value = 1106.4
$ws.Range("M55").Value = -92.66665999999998
$ws.Range("N55").Value = -1452.4

$ws.Range("H58").Value = 49999
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 49999
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 49999
$ws.Range("N58").Value = -50519

$ws.Range("H76").Value = 18000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 18000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 18000
$ws.Range("N76").Value = -18676

$ws.Range("H79").Value = 18000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 18000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 18000
$ws.Range("N79").Value = -20340

$ws.Range("H111").Value = 25000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 25000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 25000
$ws.Range("N111").Value = -33180

$ws.Range("H122").Value = 5680
$ws.Range("I122").Value = 4990.091
$ws.Range("J122").Value = 7197.8
$ws.Range("K122").Value = 14970.273
$ws.Range("L122").Value = 21593.4
$ws.Range("M122").Value = -12520.273
$ws.Range("N122").Value = -26493.4

$ws.Range("H124").Value = 22998
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 22998
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 22998
$ws.Range("N124").Value = -32818

$ws.Range("H126").Value = 5373.25
$ws.Range("I126").Value = 4996.5
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 14989.5
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -12519.5
$ws.Range("N126").Value = -22190

$ws.Range("H133").Value = 49999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059

$ws.Range("H141").Value = 49999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 49999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 49999
$ws.Range("N141").Value = -60359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H132").Value = 41305.92
$ws.Range("I132").Value = 46711.727
$ws.Range("J132").Value = 1663.3334
$ws.Range("K132").Value = 140135.181
$ws.Range("L132").Value = 4990.0002
$ws.Range("M132").Value = -137605.181
$ws.Range("N132").Value = -10050.0002

$ws.Range("H133").Value = 87499
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 87499
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 87499
$ws.Range("N133").Value = -97619

$ws.Range("H135").Value = 74715
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 74715
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 74715
$ws.Range("N135").Value = -84855

$ws.Range("H136").Value = 4109
$ws.Range("I136").Value = 2172.5
$ws.Range("J136").Value = 5400
$ws.Range("K136").Value = 6517.5
$ws.Range("L136").Value = 16200
$ws.Range("M136").Value = -3967.5
$ws.Range("N136").Value = -21300

$ws.Range("H141").Value = 98000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 98000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 98000
$ws.Range("N141").Value = -108360
